$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix the program_identifier on row 38: it referred to the 2002 variant,
#     but this row is actually the 2006 "Hartz Reform" entry, so it gets its
#     own identifier string.
$ws.Range("A38").Value = "unemploymentBenefits2006"

# --- Add the new "Negative Income Tax Experiment" program as row 42 ---
# (values are entered in the same order the source data was assembled in,
#  so new shared-string entries line up with the original edit)

# Links column (J) holds the paper URL as a real hyperlink
$ws.Range("J42").Value = "https://ideas.repec.org/p/iza/izadps/dp2067.html"

# Papers column (I) gets the dark-grey font used for this entry's citation
$ws.Range("I42").Value = "Spermann & Strotmann (2006)"
$ws.Range("I42").Font.Color = 3355443

$ws.Range("F42").Value = "Spermann & Strotmann (2006) evaluate the effects of what they call a negative income tax experiment. This is misleading because the experiment they are describing does not resemble what is typically considered a negative income tax as proposed by Milton Friedman. In the sense of Friedman, a negative income tax is a basic income that is gradually phased-out with increasing gross earnings. Instead, Spermann & Strotmann (2006) consider a experiment where long-term unemployed receive a  subsidy on top of their gross wage. The subsidy increases with the wage and can reach a maximum of up to 643 Deutschmark."
$ws.Range("F42").WrapText = $true

$ws.Range("D42").Value = "Tax Reform"
$ws.Range("A42").Value = "negativeIncomeTax"
$ws.Range("B42").Value = "Negative Income Tax Experiment"

$ws.Range("C42").Value = 2001
$ws.Range("E42").Value = 39.612
$ws.Range("G42").Value = 1278.225
$ws.Range("G42").WrapText = $true

# Turn the URL text in J42 into a real hyperlink
$ws.Hyperlinks.Add($ws.Cells.Item(42, 10), "https://ideas.repec.org/p/iza/izadps/dp2067.html")
# Reuse the same "Link" look already used by the other hyperlink cells
# instead of letting Excel mint a fresh (duplicate) style for it.
$ws.Cells.Item(41, 10).Copy()
$ws.Cells.Item(42, 10).PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 42 needs to be tall enough to show the wrapped long description
$ws.Rows.Item(42).RowHeight = 210

# Update selection / scroll position to reflect the newly added row
$ws.Range("A42").Select()
$excel.ActiveWindow.ScrollRow = 36
